$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename category columns ---
$ws.Range("B1").Value = "Convolutional"
$ws.Range("C1").Value = "Feedforward"
$ws.Range("D1").Value = "Other"
$ws.Range("E1").Value = "Sequence"
$ws.Range("F1").Value = "Attention"
$ws.Range("G1").Value = "Embedding"
$ws.Range("H1").Value = "word2vec"
$ws.Range("I1").Value = "Graph"
$ws.Range("J1").Value = "Deep Belief Network"
$ws.Range("K1").Value = "none"

# --- Data rows 2-17: update row labels, counts, and totals ---
$ws.Range("A2").Value = "Clone Detection"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 6

$ws.Range("A3").Value = "Code Synthesis"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 23

$ws.Range("A4").Value = "Program Repair"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 6

$ws.Range("A5").Value = "Software Categorization"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2

$ws.Range("A6").Value = "Software Energy Metrics"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1

$ws.Range("A7").Value = "Testing"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1

$ws.Range("A8").Value = "Vulnerability Detection"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4

$ws.Range("A9").Value = "bug localization"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2

$ws.Range("A10").Value = "code comprehension"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 14

$ws.Range("A11").Value = "code smell"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1

$ws.Range("A12").Value = "defect prediction"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 8

$ws.Range("A13").Value = "issue close time"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1

$ws.Range("A14").Value = "language model"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 12

$ws.Range("A15").Value = "language processing"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 6
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 13

$ws.Range("A16").Value = "none"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5
$ws.Range("L16").Value = 5

$ws.Range("A17").Value = "total"
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 56
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 99

# --- Remove old row 18 ("total" row), now that data occupies rows 1-17 ---
$ws.Rows(18).Delete()

Write-Host "Literature review results updated"
